$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the two trailing PCR control rows (old rows 14 & 15).
# ---------------------------------------------------------------------------
$ws.Range("A14:A15").EntireRow.Delete()

# ---------------------------------------------------------------------------
# 2. Insert a new column before the (old) rxn-volume column so the primer
#    labware-name column can be added; this shifts the old D/E -> E/F.
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).Insert()

# Match the new column-width grouping: columns C:E all share column C's width.
$cw = $ws.Columns.Item(3).ColumnWidth
$ws.Columns.Item(4).ColumnWidth = $cw
$ws.Columns.Item(5).ColumnWidth = $cw

# ---------------------------------------------------------------------------
# 3. Re-write the header row.
# ---------------------------------------------------------------------------
$ws.Cells.Item(1,1).Value = "Sample"
$ws.Cells.Item(1,2).Value = "TECAN_sample_labware_name"
$ws.Cells.Item(1,3).Value = "TECAN_sample_location"
$ws.Cells.Item(1,4).Value = "TECAN_primer_labware_name"
$ws.Cells.Item(1,5).Value = "TECAN_primer_target_position"
$ws.Cells.Item(1,6).Value = "Notes"

# ---------------------------------------------------------------------------
# 4. Re-write the sample rows (2-13) with the new columns / values.
#    Columns: A Sample, B labware name, C location, D primer labware,
#             E primer target position, F notes.
# ---------------------------------------------------------------------------
$rows = @(
  @("Poop1",      "Sample plate 1", 1,  "N7-S5_1", 1,  "DNA in 96-well plate"),
  @("Poop2",      "Sample plate 1", 2,  "N7-S5_1", 2,  "DNA in 96-well plate"),
  @("Poop3",      "Sample plate 1", 3,  "N7-S5_1", 3,  "DNA in 96-well plate"),
  @("Soil1",      "Sample plate 1", 4,  "N7-S5_1", 10, "DNA in 96-well plate"),
  @("Soil2",      "Sample plate 1", 5,  "N7-S5_1", 12, "DNA in 96-well plate"),
  @("Soil3",      "Sample plate 1", 6,  "N7-S5_1", 96, "DNA in 96-well plate"),
  @("Tissue1",    "Sample plate 1", 7,  "N7-S5_2", 24, "DNA in 96-well plate"),
  @("Tissue2",    "Sample plate 1", 8,  "N7-S5_2", 25, "DNA in 96-well plate"),
  @("Tissue3",    "Sample plate 1", 9,  "N7-S5_2", 26, "DNA in 96-well plate"),
  @("Random1",    "Sample plate 1", 10, "N7-S5_2", 27, "DNA in 96-well plate"),
  @("Random2",    "Sample plate 1", 11, "N7-S5_2", 28, "DNA in 96-well plate"),
  @("DNAx_blank", "Sample plate 1", 12, "N7-S5_2", 29, "DNA extraction blank")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r,1).Value = $row[0]
    $ws.Cells.Item($r,2).Value = $row[1]
    $ws.Cells.Item($r,3).Value = $row[2]
    $ws.Cells.Item($r,4).Value = $row[3]
    $ws.Cells.Item($r,5).Value = $row[4]
    $ws.Cells.Item($r,6).Value = $row[5]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 5. Update the selected cell to match the author's final selection.
# ---------------------------------------------------------------------------
$ws.Range("C17").Select() | Out-Null
